# edit.ps1 - apply the diff to FinalProjectPoster.pptx (single-slide deck)
#
# Summary of changes on Slide 1:
#   1. Shape 68 ("Background" box): shrink height (cy 8373300 -> 7819200 EMU)
#   2. Shape 68, RQ1 bullet: reword and split into 3 runs, middle run "Blue Bikes "
#      gets dk1 scheme-color fill
#   3. Shape 68, H1 bullet: "will be" -> "was"
#   4. Shape 68, RQ2 bullet: reword ("is"->"was", "distance"->"time", "from"->"to")
#   5. Shape 68, H2 bullet (2nd run): "from"->"to", "will be"->"was"
#   6. Shape 69 ("Primary data sources" box): reword "found median trip time..."
#      bullet
#   7. Shape 90 (regression equation box): split run "StartWalkTime" into
#      "First" + "WalkTime" (two runs, identical italic/dk1 formatting)
#
# NOTE: this interpreter's PowerShell "function" blocks do not reliably keep
# live references to the underlying COM objects, so the run-splitting logic
# (insert a unique marker character range, then delete it again in order to
# force a run boundary without changing the visible text) is inlined at each
# call site rather than factored into a reusable function.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Shape 68 ("Background") -> slide Shapes index 9 (id 60 == index 1)
# ---------------------------------------------------------------------------
$shBackground = $s.Shapes.Item(9)

$tfBackground = $shBackground.TextFrame
$trBackground = $tfBackground.TextRange

# 2) RQ1 bullet (paragraph 7) - reword + split into 3 runs
$paraRQ1 = $trBackground.Paragraphs(7, 1)
$rqRun1 = $paraRQ1.Runs(1, 1)

$rq1Part1 = "RQ1: In Cambridge, MA in summer 2023, for what percent of start/stop "
$rq1Part2 = "Blue Bikes "
$rq1Part3 = " station pairs was bikeshare faster (slower) than the T (within Cambridge)? How much faster (slower) was bikeshare?"

$rqRun1.Text = $rq1Part1 + $rq1Part2 + $rq1Part3

$marker = "@@SPLITMARK@@"

# -- split boundary between part1 and part2 --
# (the whole run 1 text is exactly part1+part2+part3, so the tail length
# from pos1 to the end of the run is part2.Length + part3.Length)
$pos1 = $rq1Part1.Length + 1
$tailLen1 = $rq1Part2.Length + $rq1Part3.Length
$tail1 = $paraRQ1.Characters($pos1, $tailLen1)
$null = $tail1.InsertBefore($marker)
$full1 = $paraRQ1.Text
$mIdx1 = $full1.IndexOf($marker)
$mSub1 = $paraRQ1.Characters($mIdx1 + 1, $marker.Length)
$mSub1.Text = ""

# -- split boundary between part2 and part3 --
$pos2 = $rq1Part1.Length + $rq1Part2.Length + 1
$tailLen2 = $rq1Part3.Length
$tail2 = $paraRQ1.Characters($pos2, $tailLen2)
$null = $tail2.InsertBefore($marker)
$full2 = $paraRQ1.Text
$mIdx2 = $full2.IndexOf($marker)
$mSub2 = $paraRQ1.Characters($mIdx2 + 1, $marker.Length)
$mSub2.Text = ""

$rqRun2 = $paraRQ1.Runs(2, 1)
$rqRun2.Font.Color.ObjectThemeColor = 1   # msoThemeColorDark1 -> schemeClr dk1

# 3) H1 bullet (paragraph 8)
$paraH1 = $trBackground.Paragraphs(8, 1)
$paraH1.Runs(1, 1).Text = "H1: Bikeshare was faster"

# 4) RQ2 bullet (paragraph 9)
$paraRQ2 = $trBackground.Paragraphs(9, 1)
$paraRQ2.Runs(1, 1).Text = "RQ2: What was the relationship between initial walking time to a T station and the difference in travel times across modes?"

# 5) H2 bullet (paragraph 10), second run holds the colored text
$paraH2 = $trBackground.Paragraphs(10, 1)
$paraH2.Runs(2, 1).Text = "As initial walking distance to a T station increases, bikeshare was relatively faster"

# 6) Resize the box: cy 8373300 -> 7819200 EMU (914400 EMU = 72 pt = 1 inch)
# Must happen AFTER all text edits above: this shape uses <a:spAutoFit/>, so
# the host recomputes/overwrites the shape height whenever its text changes.
$shBackground.Height = 7819200 / 914400 * 72

# ---------------------------------------------------------------------------
# Shape 69 ("Primary data sources") -> slide Shapes index 10
# ---------------------------------------------------------------------------
$shDataSources = $s.Shapes.Item(10)
$trDataSources = $shDataSources.TextFrame.TextRange
$paraMedian = $trDataSources.Paragraphs(21, 1)
$paraMedian.Runs(1, 1).Text = "For each start/end bikeshare station pair (with > 10 rides), calculated difference in median trip times for bikeshare vs T/walking"

# ---------------------------------------------------------------------------
# Shape 90 (regression equation "Model" box) -> slide Shapes index 31
# Split run "StartWalkTime" into "First" + "WalkTime"
# ---------------------------------------------------------------------------
$shModel = $s.Shapes.Item(31)
$trModel = $shModel.TextFrame.TextRange
$paraEq = $trModel.Paragraphs(1, 1)

# Paragraph text is: "DiffInTravelTimes = -8.24 - 0.38*(StartWalkTime)"
# "StartWalkTime" is run 3 (run 4 is the closing paren ")"); split run 3
# right after "Start" (5 chars in) -- NOT all the way to the end of the
# paragraph, or the tail range would swallow run 4 ")" as well.
$eqFullText = $paraEq.Text
$swIdx = $eqFullText.IndexOf("StartWalkTime")       # 0-based
$splitPos = $swIdx + 5 + 1                          # 1-based position of "W"
$tailLen3 = "WalkTime".Length                       # stop at end of run 3

$tail3 = $paraEq.Characters($splitPos, $tailLen3)
$null = $tail3.InsertBefore($marker)
$full3 = $paraEq.Text
$mIdx3 = $full3.IndexOf($marker)
$mSub3 = $paraEq.Characters($mIdx3 + 1, $marker.Length)
$mSub3.Text = ""

$eqRun3 = $paraEq.Runs(3, 1)
$eqRun3.Text = "First"
